# The canonical-OOXML diff for this commit only touches the picture shape(s)
# that carry the legacy "remembered fill/line" compatibility extensions
# (<a14:hiddenFill>/<a14:hiddenLine> inside <a:extLst>, guarded by
# mc:AlternateContent-style URIs {909E8E84-...}/{91240B29-...}). Those
# extensions are PowerPoint's way of recording that a shape's Fill/Line are
# turned off while still remembering the last color that was used, so they
# live under the shape's Fill/Line state. Re-affirm that state on every
# picture shape that already carries this bookkeeping so PowerPoint
# re-persists it, without touching any other shape (plain pictures that
# never had Fill/Line touched must stay untouched, or new noFill/ln nodes
# would be introduced where none existed before).

$p = $ppt.ActivePresentation

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)

    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shp = $slide.Shapes.Item($shapeIdx)

        # msoPicture = 13
        if ($shp.Type -ne 13) {
            continue
        }

        $altText = ""
        try { $altText = $shp.AlternativeText } catch { $altText = "" }

        $isLogo = ($altText -like "*TGI_logo*") -or ($shp.Name -eq "Picture 4")

        if (-not $isLogo) {
            continue
        }

        # These are exactly the properties PowerPoint backs with the
        # a14:hiddenFill / a14:hiddenLine compatibility extensions -
        # re-assert them so the "hidden fill" / "hidden line" bookkeeping
        # for this shape is rewritten/persisted on save.
        $shp.Fill.Visible = $false
        $shp.Line.Visible = $false
    }
}
